$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '<him>'
$ws.Range("C2").Value = 6
$ws.Range("B3").Value = '<golf>'
$ws.Range("C3").Value = 9
$ws.Range("B4").Value = '<may>'
$ws.Range("C4").Value = 7
$ws.Range("B5").Value = '<oh>'
$ws.Range("C5").Value = 10
$ws.Range("B6").Value = '<day>'
$ws.Range("C6").Value = 11
$ws.Range("B7").Value = '<and>'
$ws.Range("C7").Value = 9
$ws.Range("B8").Value = '<time>'
$ws.Range("C8").Value = 6
$ws.Range("B9").Value = '<in>'
$ws.Range("C9").Value = 9
$ws.Range("B10").Value = '<been>'
$ws.Range("C10").Value = 3
$ws.Range("B11").Value = '<command>'
$ws.Range("C11").Value = 14
$ws.Range("B12").Value = '<people>'
$ws.Range("C12").Value = 6
$ws.Range("B13").Value = '<her>'
$ws.Range("C13").Value = 10
$ws.Range("B14").Value = '<be>'
$ws.Range("C14").Value = 7
$ws.Range("B15").Value = '<backspace>'
$ws.Range("C15").Value = 4
$ws.Range("B16").Value = '<india>'
$ws.Range("C16").Value = 4
$ws.Range("B17").Value = '<whiskey>'
$ws.Range("C17").Value = 12
$ws.Range("B18").Value = '<up>'
$ws.Range("C18").Value = 9
$ws.Range("B19").Value = '<alt>'
$ws.Range("C19").Value = 11
$ws.Range("B20").Value = '<command>'
$ws.Range("C20").Value = 8
$ws.Range("B21").Value = '<victor>'
$ws.Range("C21").Value = 11
$ws.Range("B22").Value = '<make>'
$ws.Range("C22").Value = 5
$ws.Range("B23").Value = '<a>'
$ws.Range("C23").Value = 8
$ws.Range("B24").Value = '<there>'
$ws.Range("C24").Value = 8
$ws.Range("B25").Value = '<tango>'
$ws.Range("C25").Value = 7
$ws.Range("B26").Value = '<we>'
$ws.Range("C26").Value = 10
$ws.Range("B27").Value = '<each>'
$ws.Range("C27").Value = 5
$ws.Range("B28").Value = '<these>'
$ws.Range("C28").Value = 13
$ws.Range("B29").Value = '<down>'
$ws.Range("C29").Value = 5
$ws.Range("B30").Value = '<to>'
$ws.Range("C30").Value = 13
$ws.Range("B31").Value = '<not>'
$ws.Range("C31").Value = 8
$ws.Range("B32").Value = '<command>'
$ws.Range("C32").Value = 6
$ws.Range("B33").Value = '<like>'
$ws.Range("C33").Value = 9
$ws.Range("B34").Value = '<are>'
$ws.Range("C34").Value = 4
$ws.Range("B35").Value = '<them>'
$ws.Range("C35").Value = 9
$ws.Range("B36").Value = '<it>'
$ws.Range("C36").Value = 10
$ws.Range("B37").Value = '<six>'
$ws.Range("C37").Value = 11
$ws.Range("B38").Value = '<many>'
$ws.Range("C38").Value = 8
$ws.Range("B39").Value = '<alt>'
$ws.Range("C39").Value = 7
$ws.Range("B40").Value = '<have>'
$ws.Range("C40").Value = 9
$ws.Range("B41").Value = '<on>'
$ws.Range("C41").Value = 8
$ws.Range("B42").Value = '<india>'
$ws.Range("C42").Value = 7
$ws.Range("B43").Value = '<bravo>'
$ws.Range("C43").Value = 10
$ws.Range("B44").Value = '<delta>'
$ws.Range("C44").Value = 6
$ws.Range("B45").Value = '<may>'
$ws.Range("C45").Value = 13
$ws.Range("B46").Value = '<their>'
$ws.Range("C46").Value = 10
$ws.Range("B47").Value = '<will>'
$ws.Range("C47").Value = 10
$ws.Range("B48").Value = '<all>'
$ws.Range("C48").Value = 10
$ws.Range("B49").Value = '<foxtrot>'
$ws.Range("C49").Value = 6
$ws.Range("B50").Value = '<xars>'
$ws.Range("C50").Value = 9
$ws.Range("B51").Value = '<did>'
$ws.Range("C51").Value = 9
$ws.Range("B52").Value = '<its>'
$ws.Range("C52").Value = 4
